$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new RS Online heatsink URL as a hyperlink in the row below the last component (row 13)
$ws.Hyperlinks.Add($ws.Range("E13"), "https://nl.rs-online.com/web/p/heatsinks/1898628/")

# Match formatting of the other hyperlink cells in column E
$ws.Range("E13").Style = $ws.Range("E12").Style

# Restore the selection state recorded in the saved workbook
$ws.Range("E25").Select()
